# Insert a new data row above row 53 (shifts existing rows 53:178 down to 54:179)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("53").Insert()

# Populate the newly inserted row 53 with its data (same constant columns as
# the surrounding rows, plus the new record's own values).
$ws.Range("A53").Value = 7
$ws.Range("B53").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C53").Value = "Ñuble"
$ws.Range("D53").Value = 44498
$ws.Range("E53").Value = 16
$ws.Range("F53").Value = 100112008
$ws.Range("G53").Value = "Coliflor"
$ws.Range("H53").Value = "Sin especificar"
$ws.Range("I53").Value = "Primera"
$ws.Range("J53").Value = 600
$ws.Range("K53").Value = 700
$ws.Range("L53").Value = 750
$ws.Range("M53").Value = 725
$ws.Range("N53").Value = "$/unidad"
$ws.Range("O53").Value = "Región del Maule"
$ws.Range("P53").Value = 725
$ws.Range("Q53").Value = 1
$ws.Range("R53").Value = "Hortaliza"
